$wb = $excel.ActiveWorkbook
$lastIndex = $wb.Worksheets.Count
$sourceSheet = $wb.Worksheets.Item($lastIndex)
$sourceSheet.Copy($null, $sourceSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Tir_430_50R38"
$newSheet.Range("H3").Value = "Delft_430_50R38"
$newSheet.Range("H5").Value = "which('Truck_430_50R38.tir')"

$firstSheet = $wb.Worksheets.Item(1)
$firstSheet.Activate()
$firstSheet.Range("C27").Select()

$newSheet.Activate()
$newSheet.Range("H6").Select()
